$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Colours used by the workbook's built-in cell styles (COM uses BGR
# ordering), matched against the "Good" / "Neutral" / "Bad" cell
# styles already defined in the workbook:
#   Good    -> fill FFC6EFCE / font FF006100
#   Neutral -> fill FFFFEB9C / font FF9C5700
#   Bad     -> fill FFFFC7CE / font FF9C0006
# -----------------------------------------------------------------
$goodFill    = 13561798
$goodFont    = 24832
$neutralFill = 10284031
$neutralFont = 22428

# -----------------------------------------------------------------
# Row 18: a new background/profession entry ("Furtivité") is marked
# complete (Good) and filled in with its rolled talents; a
# specialisation ("Duelliste") is added and marked in-progress
# (Neutral).
# -----------------------------------------------------------------
$ws.Range("C18").Interior.Color = $goodFill
$ws.Range("C18").Font.Color = $goodFont

$ws.Range("E18").Value = "Malchanceux, Malentendant"
$ws.Range("E18").Interior.Color = $goodFill
$ws.Range("E18").Font.Color = $goodFont

$ws.Range("F18").Interior.Color = $goodFill
$ws.Range("F18").Font.Color = $goodFont

$ws.Range("G18").Interior.Color = $goodFill
$ws.Range("G18").Font.Color = $goodFont

$ws.Range("H18").Interior.Color = $goodFill
$ws.Range("H18").Font.Color = $goodFont

$ws.Range("J18").Value = "Duelliste"
$ws.Range("J18").Interior.Color = $neutralFill
$ws.Range("J18").Font.Color = $neutralFont

# -----------------------------------------------------------------
# Row 19 & 20: the skill entries ("Intuition", "Investigation") are
# marked complete (Good).
# -----------------------------------------------------------------
$ws.Range("C19").Interior.Color = $goodFill
$ws.Range("C19").Font.Color = $goodFont

$ws.Range("I19").Interior.Color = $goodFill
$ws.Range("I19").Font.Color = $goodFont

$ws.Range("C20").Interior.Color = $goodFill
$ws.Range("C20").Font.Color = $goodFont

# -----------------------------------------------------------------
# Column I (talent pairs, rows 19-31): the talent list shifts up by
# one entry now that row 18 claims "Apothicairerie, Toxicologie",
# dropping "Monture de combat, Obé. ànimale" and re-appending
# "Alchimie, Artillerie" at the bottom of the list.
# -----------------------------------------------------------------
$ws.Range("I19").Value = "Attaque incapacitante / sanglante"
$ws.Range("I20").Value = "Attaque magique / précise"
$ws.Range("I21").Value = "Attaques multiples, Double attaque"
$ws.Range("I22").Value = "Blindé, Protection"
$ws.Range("I23").Value = "Chasse-mage, Sentinelle"
$ws.Range("I24").Value = "Critique supérieur, Initiative supérieure"
$ws.Range("I25").Value = "Érudition, Mythologie"
$ws.Range("I26").Value = "Frappe précise, Lancer précis"
$ws.Range("I27").Value = "Interrogateur, Résilience"
$ws.Range("I28").Value = "Magie guerrière, Magie précise"
$ws.Range("I30").Value = "Onde de choc, Otage"
$ws.Range("I31").Value = "Alchimie, Artillerie"

# Leave the selection on the last touched cell, as in the source file.
$ws.Range("I31").Select()
